# Experimento2: re-sorted results table (by F1 Score, descending) with
# refreshed metric values from a re-run of the experiments, plus a bordered
# header row (matching the header style already used on the other
# "Experimento" sheets).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experimento2")

# --- Header row (A1:F1): bold + thin box border around every header cell ---
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# --- Data rows 2-36: Model name, Accuracy, AUC, Precision, Recall, F1 ---
$rows = @(
    @(2,  "RF OVERSAMPLE",             0.99775, 0.625,               1,                   0.25,                0.4),
    @(3,  "RF SMOTE",                  0.996,   0.6313277542008752,  0.7142857142857143,  0.2631578947368421,  0.3846153846153846),
    @(4,  "RF IMABALANCED",            0.99675, 0.5998745294855709,  0.75,                0.2,                 0.3157894736842106),
    @(5,  "XGBOOST IMBALANCED",        0.99675, 0.5998745294855709,  0.75,                0.2,                 0.3157894736842106),
    @(6,  "RF ADASYN",                 0.996,   0.5832077682906412,  0.75,                0.1666666666666667,  0.2727272727272728),
    @(7,  "XGBOOST OVERSAMPLE",        0.996,   0.624122367101304,   0.3,                 0.25,                0.2727272727272727),
    @(8,  "XGBOOST ADASYN ",           0.99325, 0.6371309782911992,  0.2631578947368421,  0.2777777777777778,  0.2702702702702702),
    @(9,  "XGBOOST SMOTE",             0.9925,  0.6295694020280543,  0.2380952380952381,  0.2631578947368421,  0.25),
    @(10, "KNN IMBALANCE",             0.99675, 0.5666666666666667,  1,                   0.1333333333333333,  0.2352941176470588),
    @(11, "DT IMBALANCED",             0.995,   0.5989962358845672,  0.2727272727272727,  0.2,                 0.2307692307692308),
    @(12, "DT OVERSAMPLE",             0.99575, 0.5824557004346373,  0.2222222222222222,  0.1666666666666667,  0.1904761904761905),
    @(13, "NB SMOTE",                  0.97525, 0.6994738164174565,  0.08333333333333333, 0.4210526315789473,  0.1391304347826087),
    @(14, "DT SMOTE",                  0.98525, 0.5997369082087283,  0.08333333333333333, 0.2105263157894737,  0.1194029850746269),
    @(15, "LR SMOTE",                  0.987,   0.5744258914052275,  0.07692307692307693, 0.1578947368421053,  0.103448275862069),
    @(16, "KNN OVERSAMPLE",            0.991,   0.5800735539953193,  0.07142857142857142, 0.1666666666666667,  0.1),
    @(17, "NB IMBALANCED",             0.9785,  0.6239230447511501,  0.05063291139240506, 0.2666666666666667,  0.0851063829787234),
    @(18, "NB ADASYN",                 0.96475, 0.6504687761593838,  0.04444444444444445, 0.3333333333333333,  0.07843137254901961),
    @(19, "LR OVERSAMPLE",             0.982,   0.6171013039117352,  0.04545454545454546, 0.25,                0.07692307692307693),
    @(20, "LR ADASYN ",                0.97975, 0.5750460405156538,  0.04347826086956522, 0.1666666666666667,  0.06896551724137931),
    @(21, "NB OVERSAMPLE",             0.96575, 0.6504931461049815,  0.03007518796992481, 0.3333333333333333,  0.05517241379310345),
    @(22, "MLPClassifier  ADASYN",     0.94475, 0.6404235727440146,  0.02790697674418605, 0.3333333333333333,  0.05150214592274679),
    @(23, "NB UNDERSAMPLE",            0.96125, 0.6152655792555416,  0.02702702702702703, 0.2666666666666667,  0.049079754601227),
    @(24, "DT ADASYN",                 0.9795,  0.5472682627378761,  0.02941176470588235, 0.1111111111111111,  0.04651162790697674),
    @(25, "MLPClassifier SMOTE",       0.854,   0.7171300519573236,  0.01873935264054515, 0.5789473684210527,  0.03630363036303631),
    @(26, "KNN SMOTE",                 0.73775, 0.6325374476130038,  0.009523809523809525,0.5263157894736842,  0.01870907390084191),
    @(27, "DT UNDERSAMPLE",            0.63225, 0.6493935591802592,  0.006775067750677507,0.6666666666666666,  0.01341381623071764),
    @(28, "RF UNDERSAMPLE",            0.657,   0.5289836888331242,  0.004382761139517896,0.4,                 0.008670520231213872),
    @(29, "XGBOOST UNDERSAMPLE",       0.3085,  0.5533249686323714,  0.004324324324324324,0.8,                 0.00860215053763441),
    @(30, "KNN UNDERSAMPLE",           0.53525, 0.5342952739439565,  0.004301075268817204,0.5333333333333333,  0.008533333333333334),
    @(31, "KNN ADASYN",                0.759,   0.4918243205536024,  0.00419287211740042, 0.2222222222222222,  0.00823045267489712),
    @(32, "MLPClassifier OVERSAMPLE",  0.18675, 0.5090688732865263,  0.00306654400490647, 0.8333333333333334,  0.006110601894286587),
    @(33, "MLPClassifier IMBALANCE",   0.996,   0.4998745294855709,  0,                   0,                   0),
    @(34, "MLPClassifier UNDERSAMPLE", 0.9935,  0.4986198243412798,  0,                   0,                   0),
    @(35, "LR IMBALANCED",             0.99625, 0.5,                 0,                   0,                   0),
    @(36, "LR UNDERSAMPLE",            0.99625, 0.5,                 0,                   0,                   0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
